# Change the sort-key separator from ":" to "|" on the "psami" sheet.
# Colons can appear inside data values (e.g. times in date/time sort keys),
# so "|" is used as the partition/sort key field separator instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("psami")

$ws.Range("B2").Value = "tourney|20201225"
$ws.Range("B3").Value = "joined|Tom"
$ws.Range("B4").Value = "joined|Valeria"
$ws.Range("B5").Value = "joined|Trevor"
$ws.Range("B6").Value = "game|1"
$ws.Range("B7").Value = "game|2"
$ws.Range("B8").Value = "game|3"

# Move the active selection to B9 (matches the saved cursor position).
$ws.Activate()
$ws.Range("B9").Select()
